$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: mean
$ws.Range("B3").Value = 0.9485561730412988
$ws.Range("C3").Value = 0.9436283734534584
$ws.Range("D3").Value = 0.9485561730412988
$ws.Range("E3").Value = 0.938428932084366

# Row 4: std
$ws.Range("B4").Value = 0.04968221054436373
$ws.Range("C4").Value = 0.05221303039029642
$ws.Range("D4").Value = 0.04968221054436373
$ws.Range("E4").Value = 0.05227333789916834

# Row 5: min
$ws.Range("B5").Value = 0.7744088482074752
$ws.Range("C5").Value = 0.7734731084232268
$ws.Range("D5").Value = 0.7744088482074752
$ws.Range("E5").Value = 0.7714072983942493

# Row 6: 25%
$ws.Range("B6").Value = 0.9405987795575896
$ws.Range("C6").Value = 0.9304603750798308
$ws.Range("D6").Value = 0.9405987795575896
$ws.Range("E6").Value = 0.9245407920573188

# Row 7: 50%
$ws.Range("B7").Value = 0.9610983981693364
$ws.Range("C7").Value = 0.9569284490317792
$ws.Range("D7").Value = 0.9610983981693364
$ws.Range("E7").Value = 0.9550552848161827

# Row 8: 75%
$ws.Range("B8").Value = 0.9814073226544622
$ws.Range("C8").Value = 0.9792672555876767
$ws.Range("D8").Value = 0.9814073226544622
$ws.Range("E8").Value = 0.9731010075702327
